$d = $word.ActiveDocument

# Collapse the split runs in the Title paragraph into a single run.
$find = $d.Content.Find
$find.Execute("Factsheet: Greek letters", $true, $false, $false, $false, $false, $true, 1, $false, "Factsheet: Greek letters", 2)

# Collapse the split runs in the Author paragraph into a single run.
$find = $d.Content.Find
$find.Execute("Tom Coleman", $true, $false, $false, $false, $false, $true, 1, $false, "Tom Coleman", 2)

# Collapse the split runs in the Abstract paragraph into a single run.
$find = $d.Content.Find
$find.Execute("Greek letters and their pronunciations in English.", $true, $false, $false, $false, $false, $true, 1, $false, "Greek letters and their pronunciations in English.", 2)
